# Applies the cryptos.xlsx price/volume refresh + two adjacent-row coin swaps
# described in the commit 'Updated cryptos list ... with GitHub Actions'.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Price/Volume/Coin/Link cells are stored as text (even when the text looks
# like a number, e.g. "1.00" or "0.0000254"). Assigning a plain numeric-looking
# string to .Value lets Excel auto-convert it to a Number, so we force text with a
# leading apostrophe and then clear the resulting quote-prefix style so no stray
# formatting is left behind on the cell.
function Set-TextValue {
    param($Cell, $Text)
    $Cell.Value = "'" + $Text
    $Cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "64.056.87"
Set-TextValue $ws.Range("E2") "  -3.38%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.124.50"
Set-TextValue $ws.Range("E3") "  -2.84%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.13%  "

# Row 5
Set-TextValue $ws.Range("D5") "607.66"
Set-TextValue $ws.Range("E5") "  +0.22%  "

# Row 6
Set-TextValue $ws.Range("D6") "146.95"
Set-TextValue $ws.Range("E6") "  -5.25%  "

# Row 7
Set-TextValue $ws.Range("E7") "  +0.13%  "

# Row 8
Set-TextValue $ws.Range("D8") "3.120.42"
Set-TextValue $ws.Range("E8") "  -2.94%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.524"
Set-TextValue $ws.Range("E9") "  -3.77%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.150"
Set-TextValue $ws.Range("E10") "  -5.75%  "

# Row 11
Set-TextValue $ws.Range("D11") "5.52"
Set-TextValue $ws.Range("E11") "  -3.16%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.473"
Set-TextValue $ws.Range("E12") "  -5.56%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.0000254"
Set-TextValue $ws.Range("E13") "  -4.66%  "

# Row 14
Set-TextValue $ws.Range("D14") "36.21"
Set-TextValue $ws.Range("E14") "  -5.35%  "

# Row 15
Set-TextValue $ws.Range("D15") "3.642.50"
Set-TextValue $ws.Range("E15") "  -2.63%  "

# Row 16
Set-TextValue $ws.Range("D16") "64.115.26"
Set-TextValue $ws.Range("E16") "  -3.39%  "

# Row 17
Set-TextValue $ws.Range("E17") "  -0.13%  "

# Row 18
Set-TextValue $ws.Range("D18") "3.114.22"
Set-TextValue $ws.Range("E18") "  -3.15%  "

# Row 19
Set-TextValue $ws.Range("D19") "6.89"
Set-TextValue $ws.Range("E19") "  -4.76%  "

# Row 20
Set-TextValue $ws.Range("D20") "476.50"
Set-TextValue $ws.Range("E20") "  -5.76%  "

# Row 21
Set-TextValue $ws.Range("D21") "14.44"
Set-TextValue $ws.Range("E21") "  -5.06%  "

# Row 22
Set-TextValue $ws.Range("D22") "0.701"
Set-TextValue $ws.Range("E22") "  -3.61%  "

# Row 23
Set-TextValue $ws.Range("D23") "7.66"
Set-TextValue $ws.Range("E23") "  -3.96%  "

# Row 24
Set-TextValue $ws.Range("D24") "13.60"
Set-TextValue $ws.Range("E24") "  -6.10%  "

# Row 25
Set-TextValue $ws.Range("D25") "82.82"
Set-TextValue $ws.Range("E25") "  -2.50%  "

# Row 26
Set-TextValue $ws.Range("D26") "0.999"
Set-TextValue $ws.Range("E26") "  -0.05%  "

# Row 27
Set-TextValue $ws.Range("D27") "2.91"
Set-TextValue $ws.Range("E27") "  -2.84%  "

# Row 28
Set-TextValue $ws.Range("D28") "8.38"
Set-TextValue $ws.Range("E28") "  -6.69%  "

# Row 29
Set-TextValue $ws.Range("D29") "2.21"
Set-TextValue $ws.Range("E29") "  -5.93%  "

# Row 30
Set-TextValue $ws.Range("D30") "0.121"
Set-TextValue $ws.Range("E30") "  -18.81%  "

# Row 31
Set-TextValue $ws.Range("D31") "6.73"
Set-TextValue $ws.Range("E31") "  -2.51%  "

# Row 32
Set-TextValue $ws.Range("D32") "1.00"
Set-TextValue $ws.Range("E32") "  +0.01%  "

# Row 33
Set-TextValue $ws.Range("D33") "2.69"
Set-TextValue $ws.Range("E33") "  -5.98%  "

# Row 34
Set-TextValue $ws.Range("D34") "26.17"
Set-TextValue $ws.Range("E34") "  -7.28%  "

# Row 35
Set-TextValue $ws.Range("D35") "1.10"
Set-TextValue $ws.Range("E35") "  -6.00%  "

# Row 36
Set-TextValue $ws.Range("B36") "OKB"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D36") "54.47"
Set-TextValue $ws.Range("E36") "  -1.67%  "

# Row 37
Set-TextValue $ws.Range("B37") "Filecoin"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D37") "6.02"
Set-TextValue $ws.Range("E37") "  -5.61%  "

# Row 38
Set-TextValue $ws.Range("B38") "PEPE"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws.Range("D38") "0.0₃0721"
Set-TextValue $ws.Range("E38") "  -6.31%  "

# Row 39
Set-TextValue $ws.Range("B39") "dogwifhat"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D39") "3.04"
Set-TextValue $ws.Range("E39") "  +1.12%  "

# Row 40
Set-TextValue $ws.Range("D40") "447.87"
Set-TextValue $ws.Range("E40") "  -10.20%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.0395"
Set-TextValue $ws.Range("E41") "  -5.32%  "

# Row 42
Set-TextValue $ws.Range("E42") "  -5.34%  "

# Row 43
Set-TextValue $ws.Range("D43") "8.35"
Set-TextValue $ws.Range("E43") "  -4.10%  "

# Row 44
Set-TextValue $ws.Range("D44") "2.850.73"
Set-TextValue $ws.Range("E44") "  -2.44%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.268"
Set-TextValue $ws.Range("E45") "  -8.50%  "

# Row 46
Set-TextValue $ws.Range("D46") "2.24"
Set-TextValue $ws.Range("E46") "  -7.93%  "

# Row 47
Set-TextValue $ws.Range("B47") "USDe"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D47") "0.998"
Set-TextValue $ws.Range("E47") "  -0.11%  "

# Row 48
Set-TextValue $ws.Range("B48") "InjectiveProtocol"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D48") "26.28"
Set-TextValue $ws.Range("E48") "  -5.98%  "

# Row 49
Set-TextValue $ws.Range("B49") "ThetaToken"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Range("D49") "2.30"
Set-TextValue $ws.Range("E49") "  -3.62%  "

# Row 50
Set-TextValue $ws.Range("B50") "Stellar"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D50") "0.114"
Set-TextValue $ws.Range("E50") "  -3.01%  "

# Row 51
Set-TextValue $ws.Range("D51") "118.49"
Set-TextValue $ws.Range("E51") "  -2.20%  "
